{"js": "// The document originally ends with:\n//   ... \"Sfdsdsgdfgdfgdfgdfgdfg\" (para 4)\n//   \"sfdsdsgdfgdfgdfgdfgdfg\" (para 5, contains the _GoBack bookmark)\n//\n// The edit inserts 7 new paragraphs between those two, reusing the same\n// run/paragraph formatting (rFonts hint=\"default\", lang=\"en-US\") that the\n// surrounding paragraphs already use. Cloning via Paragraph.insertParagraph()\n// on the existing \"Sfdsdsgdfgdfgdfgdfgdfg\" paragraph carries that formatting\n// forward automatically, so each new paragraph doesn't need explicit font\n// settings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is exactly \"Sfdsdsgdfgdfgdfgdfgdfg\"\n// (capital S) - that's the anchor after which the new content is inserted.\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Sfdsdsgdfgdfgdfgdfgdfg\") {\n    anchor = p;\n  }\n}\n\nconst newTexts = [\n  \"Sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfgsfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n];\n\nlet current = anchor;\nfor (const text of newTexts) {\n  current = current.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# The document originally ends with:\n#   ... \"Sfdsdsgdfgdfgdfgdfgdfg\" (paragraph 4)\n#   \"sfdsdsgdfgdfgdfgdfgdfg\" (paragraph 5, holds the _GoBack bookmark)\n#\n# This inserts 7 new paragraphs between those two. Using\n# Range.InsertParagraphAfter() on an existing paragraph (and then writing the\n# new paragraph's Range.Text) clones the surrounding run/paragraph formatting\n# (rFonts hint=\"default\", lang=\"en-US\") automatically, so no explicit font\n# work is required.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is exactly \"Sfdsdsgdfgdfgdfgdfgdfg\"\n# (capital S) - comparing with .Equals() because -eq is case-insensitive here.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n  if ($t.Equals(\"Sfdsdsgdfgdfgdfgdfgdfg\")) {\n    $anchorIndex = $i\n  }\n}\n\n$newTexts = @(\n  \"Sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfgsfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\",\n  \"sfdsdsgdfgdfgdfgdfgdfg\"\n)\n\n$insertAt = $anchorIndex\nforeach ($txt in $newTexts) {\n  $p = $d.Paragraphs.Item($insertAt)\n  $p.Range.InsertParagraphAfter()\n  $insertAt = $insertAt + 1\n  $newP = $d.Paragraphs.Item($insertAt)\n  $newP.Range.Text = $txt\n}\n"}
